# Connection between boxes and Products - Test
# Append two new rows of data (row 11 and row 12) to the Products sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "A4846885A"
$ws.Cells.Item(11, 3).Value = "2025-06-09 15:07:32"

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "'48836138"
$ws.Cells.Item(12, 3).Value = "2025-06-09 15:07:33"
